# Update column F ("dSF") values on Sheet1 to reflect the repulled/recalculated
# data, per the commit "repull data, push all data, mean calculation".
# Only column F values change; all other cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 0
    8  = 2
    9  = 5
    10 = 1
    11 = 7
    12 = -11
    13 = -2
    14 = -3
    15 = -1
    16 = -1
    17 = -1
    18 = 11
    19 = -3
    20 = -4
    22 = -1
    23 = 3
    24 = 10
    25 = 1
    26 = 1
    27 = 1
    28 = 2
    31 = 2
    32 = 2
    33 = 3
    34 = 3
    39 = -8
    40 = -5
    41 = 10
    43 = -4
    44 = 2
    45 = 1
    46 = 1
    47 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
